# Insert a new data row before row 423 (pushing existing rows 423:525 down to 424:526)
# and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 423, shifting rows 423:525 down to 424:526.
$ws.Rows.Item(423).Insert()

# Populate the newly inserted row 423 with the new weekly price record.
$ws.Cells.Item(423, 1).Value = 11
$ws.Cells.Item(423, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(423, 3).Value = "Bíobío"
$ws.Cells.Item(423, 4).Value = 45275
$ws.Cells.Item(423, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(423, 5).Value = 8
$ws.Cells.Item(423, 6).Value = 100112008
$ws.Cells.Item(423, 7).Value = "Coliflor"
$ws.Cells.Item(423, 8).Value = "Sin especificar"
$ws.Cells.Item(423, 9).Value = "Primera"
$ws.Cells.Item(423, 10).Value = 1500
$ws.Cells.Item(423, 11).Value = 1000
$ws.Cells.Item(423, 12).Value = 1000
$ws.Cells.Item(423, 13).Value = 1000
$ws.Cells.Item(423, 14).Value = "`$/unidad"
$ws.Cells.Item(423, 15).Value = "Región Metropolitana"
$ws.Cells.Item(423, 16).Value = 1000
$ws.Cells.Item(423, 17).Value = 1
$ws.Cells.Item(423, 18).Value = "Hortaliza"
